$wb = $excel.ActiveWorkbook

# --- "papers" sheet: add new row 8 data (new publication) ---
$ws = $wb.Worksheets.Item("papers")

$ws.Range("B8").Value = "S. Rudra, S. Sarker, . M. Kim"
$ws.Range("C8").Value = "Simulation and electrochemical impedance spectroscopy of dye-sensitized solar cells"
$ws.Range("D8").Value = "Journal of Industrial and Engineering Chemistry"
$ws.Range("E8").Value = 2021
$ws.Range("H8").Value = "https://doi.org/10.1016/j.jiec.2021.03.010"
$ws.Range("F8").Value = "XX"
$ws.Range("G8").Value = "XX"

$ws.Range("L8").Formula = '="<li>" & B8 & "; " & C8 & "; <em>" & PROPER(D8) & "</em>, " & E8 & ", " & F8 & ", " & G8 & " (<a href=""" & H8 & """ target=""_blank"" >" & "DOI: " &  H8 & "</a>).</li>"'

$ws.Hyperlinks.Add($ws.Range("H8"), "https://doi.org/10.1016/j.jiec.2021.03.010", "", "Persistent link using digital object identifier", "")

# Update the view so the new row is visible / selected like the saved file
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("L8").Select()

# --- "recent papers" sheet: point the "most recent" slot at the new paper ---
$ws2 = $wb.Worksheets.Item("recent papers")

$ws2.Range("A1").Value = "https://doi.org/10.1016/j.jiec.2021.03.010"
$ws2.Range("A6").Value = "Simulation and electrochemical impedance spectroscopy of dye-sensitized solar cells"

$ws2.Hyperlinks.Add($ws2.Range("A1"), "https://doi.org/10.1016/j.jiec.2021.03.010", "", "Persistent link using digital object identifier", "")

$ws2.Range("A11").Select()
